# saco acentos de los TCs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the credentials / claim number used on the "Mediación" row (row 3)
$ws.Range("D3").Value = "mbalducci"
$ws.Range("E3").Value = "gw"

# Claim number must stay text (keep the leading zero), like the existing F2 cell
$ws.Range("F3").Value = "'0420194406696"

# Move the active selection from G4 to D3
$ws.Range("D3").Select()
